$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(2,3,0.8029023333333333,2.408707,0.06206726394886004,0.06206726394886004,3,4.665097333333334,13.995292,0.06181862018736071,0.06181862018736071,3.745617534160445,33.710557807444,0.003836912616123245,0.003836912616123245),
    @(3,3,0.8029023333333333,2.408707,0.06206726394886004,0.06206726394886004,3,16.41004366666667,49.230131,0.2174544675497312,0.2174544675497312,13.17566235006855,118.580961150617,0.01349680383426798,0.01349680383426798),
    @(4,3,0.8029023333333333,2.408707,0.06206726394886004,0.06206726394886004,3,17.65804266666667,52.974128,0.2339920809504103,0.2339920809504103,14.17768365916622,127.599152932496,0.01452324825029214,0.01452324825029214),
    @(5,3,0.8029023333333333,2.408707,0.06206726394886004,0.06206726394886004,3,13.44741733333333,40.342252,0.1781958071250527,0.1781958071250527,10.79696275424044,97.172664788164,0.0110601261954108,0.0110601261954108),
    @(6,3,0.8029023333333333,2.408707,0.06206726394886004,0.06206726394886004,3,5.111344666666667,15.334034,0.06773197899594202,0.06773197899594202,4.103910559337555,36.935195034038,0.004203938618119778,0.004203938618119778),
    @(7,3,0.8029023333333333,2.408707,0.06206726394886004,0.06206726394886004,3,18.172329,54.516987,0.2408070451915031,0.2408070451915031,14.590605356201,131.315448205809,0.01494623443464609,0.01494623443464609),
    @(8,3,2.159929333333333,6.479788,0.1669703754456877,0.1669703754456877,3,4.665097333333334,13.995292,0.06181862018736071,0.06181862018736071,10.07628057312178,90.686525158096,0.01032187822221798,0.01032187822221798),
    @(9,3,2.159929333333333,6.479788,0.1669703754456877,0.1669703754456877,3,16.41004366666667,49.230131,0.2174544675497312,0.2174544675497312,35.44453467691422,319.000812092228,0.03630845408912071,0.03630845408912071),
    @(10,3,2.159929333333333,6.479788,0.1669703754456877,0.1669703754456877,3,17.65804266666667,52.974128,0.2339920809504103,0.2339920809504103,38.14012432498489,343.261118924864,0.03906974560760775,0.03906974560760775),
    @(11,3,2.159929333333333,6.479788,0.1669703754456877,0.1669703754456877,3,13.44741733333333,40.342252,0.1781958071250527,0.1781958071250527,29.04547115584178,261.409240402576,0.02975342081851739,0.02975342081851739),
    @(12,3,2.159929333333333,6.479788,0.1669703754456877,0.1669703754456877,3,5.111344666666667,15.334034,0.06773197899594202,0.06773197899594202,11.04014327831022,99.36128950479201,0.01130923396263187,0.01130923396263187),
    @(13,3,2.159929333333333,6.479788,0.1669703754456877,0.1669703754456877,3,18.172329,54.516987,0.2408070451915031,0.2408070451915031,39.250946462084,353.258518158756,0.04020764274559194,0.04020764274559194),
    @(14,3,9.973171666666666,29.919515,0.7709623606054523,0.7709623606054523,3,4.665097333333334,13.995292,0.06181862018736071,0.06181862018736071,46.52581654704222,418.73234892338,0.04765982934901949,0.04765982934901948),
    @(15,3,9.973171666666666,29.919515,0.7709623606054523,0.7709623606054523,3,16.41004366666667,49.230131,0.2174544675497312,0.2174544675497312,163.6601825451628,1472.941642906465,0.1676492096263425,0.1676492096263425),
    @(16,3,9.973171666666666,29.919515,0.7709623606054523,0.7709623606054523,3,17.65804266666667,52.974128,0.2339920809504103,0.2339920809504103,176.1066908119911,1584.96021730792,0.1803990870925104,0.1803990870925104),
    @(17,3,9.973171666666666,29.919515,0.7709623606054523,0.7709623606054523,3,13.44741733333333,40.342252,0.1781958071250527,0.1781958071250527,134.1134015386422,1207.02061384778,0.1373822601111245,0.1373822601111245),
    @(18,3,9.973171666666666,29.919515,0.7709623606054523,0.7709623606054523,3,5.111344666666667,15.334034,0.06773197899594202,0.06773197899594202,50.97631780816777,458.78686027351,0.05221880641519037,0.05221880641519037),
    @(19,3,9.973171666666666,29.919515,0.7709623606054523,0.7709623606054523,3,18.172329,54.516987,0.2408070451915031,0.2408070451915031,181.235756700145,1631.121810301305,0.1856531680112651,0.1856531680112651)
)

foreach ($entry in $data) {
    $r = $entry[0]
    $ws.Range("E" + $r).Value = $entry[1]
    $ws.Range("G" + $r).Value = $entry[2]
    $ws.Range("H" + $r).Value = $entry[3]
    $ws.Range("I" + $r).Value = $entry[4]
    $ws.Range("J" + $r).Value = $entry[5]
    $ws.Range("K" + $r).Value = $entry[6]
    $ws.Range("M" + $r).Value = $entry[7]
    $ws.Range("N" + $r).Value = $entry[8]
    $ws.Range("O" + $r).Value = $entry[9]
    $ws.Range("P" + $r).Value = $entry[10]
    $ws.Range("Q" + $r).Value = $entry[11]
    $ws.Range("R" + $r).Value = $entry[12]
    $ws.Range("S" + $r).Value = $entry[13]
    $ws.Range("T" + $r).Value = $entry[14]
}
